# Refresh cryptos list: update Price (D) and Volume(1h) (E) columns
# with the latest scraped values (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.930.43"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.764.74"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "643.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "3.763.29"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "4.398.96"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "3.771.08"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "68.915.96"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "472.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("E24").Value = "  -5.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "3.913.76"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  +16.62%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "3.720.78"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("E40").Value = "  -5.38%  "
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.957"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "155.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.79%  "
